$d = $word.ActiveDocument

# Number of paragraphs before the edit (so we can locate the new ones afterwards).
$countBefore = $d.Paragraphs.Count

# The document body always ends with an unremovable paragraph mark, so the
# last "character position" of the story range sits *inside* that final
# paragraph (just before its mark). Inserting text there appends to that
# paragraph; the first inserted paragraph mark closes it off unchanged and
# the remaining marks create the new paragraphs we want. We therefore need
# one extra leading paragraph mark to end up with two clean blank
# paragraphs, the heading paragraph, and two more blank paragraphs after
# the existing final paragraph.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertAfter("`r`r`rExercice 3 : Sudoku`r`r")

# Apply the Heading2 style to the newly inserted heading paragraph (the
# third paragraph added, i.e. after the two leading blank paragraphs).
$headingPara = $d.Paragraphs.Item($countBefore + 3)
$headingPara.Style = "Heading2"
